# Update "想去人数" (F) and a couple of "Cover" (I) values across the
# four sheets of the workbook, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 442
$ws.Range("F6").Value  = 1291
$ws.Range("F8").Value  = 7602
$ws.Range("F11").Value = 2083
$ws.Range("F12").Value = 8239
$ws.Range("F16").Value = 5633
$ws.Range("F18").Value = 2587
$ws.Range("F19").Value = 1118
$ws.Range("F21").Value = 338
$ws.Range("F22").Value = 398
$ws.Range("F25").Value = 514
$ws.Range("F26").Value = 3400
$ws.Range("F27").Value = 36
$ws.Range("F29").Value = 2896
$ws.Range("F31").Value = 332
$ws.Range("F34").Value = 132
$ws.Range("F35").Value = 643
$ws.Range("F37").Value = 873
$ws.Range("F38").Value = 1651
$ws.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202403/BIvjhmZq1709792042233.jpeg"
$ws.Range("F39").Value = 43
$ws.Range("F41").Value = 10
$ws.Range("F42").Value = 2673
$ws.Range("F44").Value = 2276

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 44
$ws.Range("F9").Value = 112

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1315

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 1315
$ws.Range("F6").Value  = 1291
$ws.Range("F7").Value  = 7602
$ws.Range("F10").Value = 2083
$ws.Range("F11").Value = 8239
$ws.Range("F15").Value = 5633
$ws.Range("F17").Value = 2587
$ws.Range("F18").Value = 1118
$ws.Range("F20").Value = 398
$ws.Range("F25").Value = 514
$ws.Range("F26").Value = 3400
$ws.Range("F27").Value = 36
$ws.Range("F29").Value = 2896
$ws.Range("F30").Value = 332
$ws.Range("F33").Value = 44
$ws.Range("F34").Value = 643
$ws.Range("F37").Value = 873
$ws.Range("F39").Value = 1651
$ws.Range("I39").Value = "//i0.hdslb.com/bfs/openplatform/202403/BIvjhmZq1709792042233.jpeg"
$ws.Range("F40").Value = 43
$ws.Range("F42").Value = 10
$ws.Range("F43").Value = 2673
$ws.Range("F46").Value = 2276
$ws.Range("F49").Value = 112
